# Add a "Sex"/"Male" enum column before the existing "Name" column.
# This shifts the former Name/Sum/Date columns (B/C/D) one to the right
# (C/D/E), matching the DataSourceDynamicPanel rendering a new
# enum-typed property ("Sex") ahead of the other primitive-typed ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; this shifts the old C (Name values) and D
# (Sum/Date values) columns to D and E respectively, carrying their
# existing styles/number-formats along for the ride.
$ws.Columns("C:C").Insert()

# Header row
$ws.Range("B2").Value = "Sex"
$ws.Range("C2").Value = "Name"

# Data row
$ws.Range("B3").Value = "Male"
$ws.Range("C3").Value = "Test"
